$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.51648072919072
$ws.Range("C2").Value = 6.145569675762306
$ws.Range("D2").Value = 4.836593880964283
$ws.Range("E2").Value = 16.45662231652857
$ws.Range("F2").Value = 24.34937062023393
$ws.Range("K2").Value = 8.667571513691255
$ws.Range("O2").Value = 21.75015368155667

$ws.Range("B3").Value = 9.17027533169408
$ws.Range("C3").Value = 5.978901922792356
$ws.Range("D3").Value = 4.795374274760977
$ws.Range("E3").Value = 15.52626842187457
$ws.Range("F3").Value = 24.35783602768995
$ws.Range("K3").Value = 8.417167355826217
$ws.Range("O3").Value = 21.81322534259327

$ws.Range("B4").Value = 8.952472728208928
$ws.Range("C4").Value = 5.873289745481984
$ws.Range("D4").Value = 4.769532975394966
$ws.Range("E4").Value = 14.93045656831966
$ws.Range("F4").Value = 24.37059821565621
$ws.Range("K4").Value = 8.260876708963771
$ws.Range("O4").Value = 21.85763764166031

$ws.Range("B5").Value = 8.862542160691657
$ws.Range("C5").Value = 5.829467662113198
$ws.Range("D5").Value = 4.758873257366198
$ws.Range("E5").Value = 14.68175075304581
$ws.Range("F5").Value = 24.37769639679664
$ws.Range("K5").Value = 8.196645470465526
$ws.Range("O5").Value = 21.87715984649077

$ws.Range("B6").Value = 8.847542731773315
$ws.Range("C6").Value = 5.822144897733559
$ws.Range("D6").Value = 4.757095554917283
$ws.Range("E6").Value = 14.64010549543998
$ws.Range("F6").Value = 24.37898948570619
$ws.Range("K6").Value = 8.185950188363185
$ws.Range("O6").Value = 21.88048730981578

$ws.Range("B7").Value = 8.951264450970392
$ws.Range("C7").Value = 5.872701865985945
$ws.Range("D7").Value = 4.769389731554663
$ws.Range("E7").Value = 14.92712595775057
$ws.Range("F7").Value = 24.37068626856873
$ws.Range("K7").Value = 8.260012523168491
$ws.Range("O7").Value = 21.85789516805773

$ws.Range("B8").Value = 9.398275280444068
$ws.Range("C8").Value = 6.088805846639046
$ws.Range("D8").Value = 4.822494617429653
$ws.Range("E8").Value = 16.141082836171
$ws.Range("F8").Value = 24.35071711371996
$ws.Range("K8").Value = 8.581816195715877
$ws.Range("O8").Value = 21.77071714308919

$ws.Range("B9").Value = 10.22784338344551
$ws.Range("C9").Value = 6.48498629271159
$ws.Range("D9").Value = 4.922209317477546
$ws.Range("E9").Value = 18.38985695548248
$ws.Range("F9").Value = 24.37172484707983
$ws.Range("K9").Value = 9.188751054253776
$ws.Range("O9").Value = 21.64513309329778

$ws.Range("B10").Value = 10.80220830364153
$ws.Range("C10").Value = 6.75732003889713
$ws.Range("D10").Value = 4.992511965249267
$ws.Range("E10").Value = 20.02410079230791
$ws.Range("F10").Value = 24.42395079512913
$ws.Range("K10").Value = 9.643766359580628
$ws.Range("O10").Value = 21.58085853345176

$ws.Range("B11").Value = 11.05475985586905
$ws.Range("C11").Value = 6.876796906513307
$ws.Range("D11").Value = 5.023798682573288
$ws.Range("E11").Value = 20.72552826646207
$ws.Range("F11").Value = 24.45569190830027
$ws.Range("K11").Value = 9.886044951325806
$ws.Range("O11").Value = 21.5577610032522

$ws.Range("B12").Value = 11.14906124350819
$ws.Range("C12").Value = 6.921380482161253
$ws.Range("D12").Value = 5.035541896473295
$ws.Range("E12").Value = 20.98513751507764
$ws.Range("F12").Value = 24.46885627261076
$ws.Range("K12").Value = 9.976025520641688
$ws.Range("O12").Value = 21.54990236830036

$ws.Range("B13").Value = 11.12881225169541
$ws.Range("C13").Value = 6.911808333974815
$ws.Range("D13").Value = 5.03301749902724
$ws.Range("E13").Value = 20.9294926158722
$ws.Range("F13").Value = 24.46597024946929
$ws.Range("K13").Value = 9.956725541476937
$ws.Range("O13").Value = 21.55155530436764

$ws.Range("B14").Value = 11.06254524801405
$ws.Range("C14").Value = 6.880478183748684
$ws.Range("D14").Value = 5.024766921545897
$ws.Range("E14").Value = 20.74700665777481
$ws.Range("F14").Value = 24.4567520314513
$ws.Range("K14").Value = 9.893483242851289
$ws.Range("O14").Value = 21.5570966477958

$ws.Range("B15").Value = 11.02177884345361
$ws.Range("C15").Value = 6.861200929279645
$ws.Range("D15").Value = 5.01969947739841
$ws.Range("E15").Value = 20.63444754702559
$ws.Range("F15").Value = 24.45125455560755
$ws.Range("K15").Value = 9.854514801772382
$ws.Range("O15").Value = 21.56060664662012

$ws.Range("B16").Value = 10.78552080502791
$ws.Range("C16").Value = 6.74942096482035
$ws.Range("D16").Value = 4.990452942181032
$ws.Range("E16").Value = 19.97741858707068
$ws.Range("F16").Value = 24.42203683650296
$ws.Range("K16").Value = 9.627687477329436
$ws.Range("O16").Value = 21.5824920080865

$ws.Range("B17").Value = 10.63828866698727
$ws.Range("C17").Value = 6.679699166830932
$ws.Range("D17").Value = 4.972330040193686
$ws.Range("E17").Value = 19.56361952687138
$ws.Range("F17").Value = 24.40615556839093
$ws.Range("K17").Value = 9.492996367106688
$ws.Range("O17").Value = 21.59749443083587

$ws.Range("B18").Value = 10.55278697203531
$ws.Range("C18").Value = 6.639183387809893
$ws.Range("D18").Value = 4.961841139402024
$ws.Range("E18").Value = 19.32165820945405
$ws.Range("F18").Value = 24.39777280098755
$ws.Range("K18").Value = 9.429424487323349
$ws.Range("O18").Value = 21.60670121248881

$ws.Range("B19").Value = 10.52369971070441
$ws.Range("C19").Value = 6.625395223770196
$ws.Range("D19").Value = 4.958278739483635
$ws.Range("E19").Value = 19.23905396751062
$ws.Range("F19").Value = 24.39506372327378
$ws.Range("K19").Value = 9.407820200038328
$ws.Range("O19").Value = 21.60991755345162

$ws.Range("B20").Value = 10.65404700959027
$ws.Range("C20").Value = 6.687164176477673
$ws.Range("D20").Value = 4.9742660235742
$ws.Range("E20").Value = 19.60807823426329
$ws.Range("F20").Value = 24.40776837536794
$ws.Range("K20").Value = 9.504723621468585
$ws.Range("O20").Value = 21.59583756276146

$ws.Range("B21").Value = 11.08204626844021
$ws.Range("C21").Value = 6.88969870471875
$ws.Range("D21").Value = 5.027193184833818
$ws.Range("E21").Value = 20.80076994803954
$ws.Range("F21").Value = 24.45942861275861
$ws.Range("K21").Value = 9.912107166817583
$ws.Range("O21").Value = 21.55544488738953

$ws.Range("B22").Value = 11.35396088979236
$ws.Range("C22").Value = 7.018211269704115
$ws.Range("D22").Value = 5.061173671832623
$ws.Range("E22").Value = 21.54526994555979
$ws.Range("F22").Value = 24.49986109966872
$ws.Range("K22").Value = 10.17069553988705
$ws.Range("O22").Value = 21.5342232299846

$ws.Range("B23").Value = 11.20957287717492
$ws.Range("C23").Value = 6.949982286850164
$ws.Range("D23").Value = 5.043094976288632
$ws.Range("E23").Value = 21.15110724077332
$ws.Range("F23").Value = 24.47767274343278
$ws.Range("K23").Value = 10.03363346111813
$ws.Range("O23").Value = 21.54507442634607

$ws.Range("B24").Value = 10.6469253229613
$ws.Range("C24").Value = 6.683790588622301
$ws.Range("D24").Value = 4.97339098253647
$ws.Range("E24").Value = 19.58799111030701
$ws.Range("F24").Value = 24.40703689632963
$ws.Range("K24").Value = 9.499423301351637
$ws.Range("O24").Value = 21.59658482115928

$ws.Range("B25").Value = 10.0091726872347
$ws.Range("C25").Value = 6.380973029697469
$ws.Range("D25").Value = 4.895734711069881
$ws.Range("E25").Value = 17.75155183517432
$ws.Range("F25").Value = 24.35958298269789
$ws.Range("K25").Value = 9.027665718823757
$ws.Range("O25").Value = 21.67421395518986
